$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range('D2') '28.702.77'
Set-TextValue $ws.Range('E2') '  +1.77%  '

Set-TextValue $ws.Range('D3') '1.809.53'
Set-TextValue $ws.Range('E3') '  -0.47%  '

Set-TextValue $ws.Range('D4') '1.002'
Set-TextValue $ws.Range('E4') '  +0.26%  '

Set-TextValue $ws.Range('D5') '328.40'
Set-TextValue $ws.Range('E5') '  -3.03%  '

Set-TextValue $ws.Range('D6') '0.9999'
Set-TextValue $ws.Range('E6') '  +0.24%  '

Set-TextValue $ws.Range('D7') '0.4378'
Set-TextValue $ws.Range('E7') '  -0.33%  '

Set-TextValue $ws.Range('D8') '0.3775'
Set-TextValue $ws.Range('E8') '  +6.74%  '

Set-TextValue $ws.Range('D9') '44.57'
Set-TextValue $ws.Range('E9') '  -2.47%  '

Set-TextValue $ws.Range('D10') '0.07702'
Set-TextValue $ws.Range('E10') '  +3.28%  '

Set-TextValue $ws.Range('D11') '1.142'
Set-TextValue $ws.Range('E11') '  -1.29%  '

Set-TextValue $ws.Range('D12') '22.72'
Set-TextValue $ws.Range('E12') '  -1.38%  '

Set-TextValue $ws.Range('D13') '1.001'
Set-TextValue $ws.Range('E13') '  +0.19%  '

Set-TextValue $ws.Range('D14') '6.288'
Set-TextValue $ws.Range('E14') '  +0.11%  '

Set-TextValue $ws.Range('D15') '7.570'
Set-TextValue $ws.Range('E15') '  +3.62%  '

Set-TextValue $ws.Range('D16') '1.806.31'
Set-TextValue $ws.Range('E16') '  -0.79%  '

Set-TextValue $ws.Range('D17') '0.00001097'
Set-TextValue $ws.Range('E17') '  +0.79%  '

Set-TextValue $ws.Range('D18') '0.06735'
Set-TextValue $ws.Range('E18') '  +0.70%  '

Set-TextValue $ws.Range('D19') '81.07'
Set-TextValue $ws.Range('E19') '  -1.27%  '

Set-TextValue $ws.Range('E20') '  +0.24%  '

Set-TextValue $ws.Range('D21') '17.72'
Set-TextValue $ws.Range('E21') '  +2.18%  '

Set-TextValue $ws.Range('D22') '6.293'
Set-TextValue $ws.Range('E22') '  -2.65%  '

Set-TextValue $ws.Range('D23') '28.697.92'
Set-TextValue $ws.Range('E23') '  +1.72%  '

Set-TextValue $ws.Range('D24') '11.82'
Set-TextValue $ws.Range('E24') '  -2.74%  '

Set-TextValue $ws.Range('D25') '2.451'
Set-TextValue $ws.Range('E25') '  +2.90%  '

Set-TextValue $ws.Range('D26') '20.60'
Set-TextValue $ws.Range('E26') '  -0.98%  '

Set-TextValue $ws.Range('D27') '154.05'
Set-TextValue $ws.Range('E27') '  -0.95%  '

Set-TextValue $ws.Range('D28') '2.376'
Set-TextValue $ws.Range('E28') '  -4.58%  '

Set-TextValue $ws.Range('D29') '2.013.50'
Set-TextValue $ws.Range('E29') '  -0.66%  '

Set-TextValue $ws.Range('D30') '1.303'
Set-TextValue $ws.Range('E30') '  -0.46%  '

Set-TextValue $ws.Range('D31') '131.58'
Set-TextValue $ws.Range('E31') '  -1.21%  '

Set-TextValue $ws.Range('D32') '3.973'
Set-TextValue $ws.Range('E32') '  -2.30%  '

Set-TextValue $ws.Range('D33') '5.838'
Set-TextValue $ws.Range('E33') '  -2.53%  '

Set-TextValue $ws.Range('D34') '0.09197'
Set-TextValue $ws.Range('E34') '  -2.36%  '

Set-TextValue $ws.Range('D35') '0.2241'
Set-TextValue $ws.Range('E35') '  +3.31%  '

Set-TextValue $ws.Range('D36') '12.24'
Set-TextValue $ws.Range('E36') '  -1.05%  '

Set-TextValue $ws.Range('D37') '0.06358'
Set-TextValue $ws.Range('E37') '  +1.40%  '

Set-TextValue $ws.Range('B38') 'InternetComputer(DFINITY)'
Set-TextValue $ws.Range('C38') 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextValue $ws.Range('D38') '5.242'
Set-TextValue $ws.Range('E38') '  -0.03%  '

Set-TextValue $ws.Range('B39') 'TheSandbox'
Set-TextValue $ws.Range('C39') 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
Set-TextValue $ws.Range('D39') '0.6625'
Set-TextValue $ws.Range('E39') '  -2.59%  '

Set-TextValue $ws.Range('B40') 'VeChain'
Set-TextValue $ws.Range('C40') 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue $ws.Range('D40') '0.02322'
Set-TextValue $ws.Range('E40') '  -2.14%  '

Set-TextValue $ws.Range('D41') '1.205'
Set-TextValue $ws.Range('E41') '  -1.15%  '

Set-TextValue $ws.Range('B42') 'FraxShare'
Set-TextValue $ws.Range('C42') 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextValue $ws.Range('D42') '8.096'
Set-TextValue $ws.Range('E42') '  -2.43%  '

Set-TextValue $ws.Range('B43') 'WEMIXTOKEN'
Set-TextValue $ws.Range('C43') 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
Set-TextValue $ws.Range('D43') '1.436'
Set-TextValue $ws.Range('E43') '  -3.41%  '

Set-TextValue $ws.Range('D44') '1.000'
Set-TextValue $ws.Range('E44') '  +0.27%  '

Set-TextValue $ws.Range('D45') '13.87'
Set-TextValue $ws.Range('E45') '  -1.33%  '

Set-TextValue $ws.Range('D46') '0.6089'
Set-TextValue $ws.Range('E46') '  -1.36%  '

Set-TextValue $ws.Range('D47') '3.801'
Set-TextValue $ws.Range('E47') '  -1.97%  '

Set-TextValue $ws.Range('D48') '128.44'
Set-TextValue $ws.Range('E48') '  -0.92%  '

Set-TextValue $ws.Range('D49') '2.034'
Set-TextValue $ws.Range('E49') '  -0.71%  '

Set-TextValue $ws.Range('D50') '0.07088'
Set-TextValue $ws.Range('E50') '  -0.39%  '

Set-TextValue $ws.Range('B51') 'EOS'
Set-TextValue $ws.Range('C51') 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
Set-TextValue $ws.Range('D51') '1.147'
Set-TextValue $ws.Range('E51') '  -2.48%  '
